$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Request" header in column F with "DOD Id"
$ws.Range("F1").Value = "DOD Id"

# Widen column F to fit the new header text
$ws.Range("F1").ColumnWidth = 11.17

# Move the selection to column G (matches the post-edit view state)
$ws.Range("G1:G1048576").Select()
